# Auto-generated PowerShell Excel COM-interop script
# Updates the cryptos list data in the active worksheet to match the target diff.
# All target cells are plain text (inlineStr) in the source workbook, so we force
# the Text number format before assignment to prevent Excel from auto-converting
# numeric-looking strings (e.g. '0.999', '311.80') into actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.906.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.479.80'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.80'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.28'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.554'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.14%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.512'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.08'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.109'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.857.56'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.477.79'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.43%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.99%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.868.58'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.42'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -5.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0922'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.79'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.05'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.32%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.03%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.63%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.77'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.64'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '154.75'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.64%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.62'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0760'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.50'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -11.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.03'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.96%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.03%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.107'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.44%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.05'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.66%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.18%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.004.04'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.09%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.86%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.70'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.716.86'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '77.64'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.83%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.183'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.33%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.07'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.62%  '
